# Update the tri proximity (5-mile / 10-mile) columns in both the
# "Means" and "Standard Deviations" summary tables.

$wb = $excel.ActiveWorkbook

# --- Means sheet ---
$wsMeans = $wb.Worksheets.Item("Means")

$wsMeans.Range("F2").Value = 79
$wsMeans.Range("G2").Value = 79

$wsMeans.Range("F3").Value = 0.25
$wsMeans.Range("G3").Value = 3

$wsMeans.Range("G4").Value = 18

$wsMeans.Range("F5").Value = 6.9
$wsMeans.Range("G5").Value = 18

$wsMeans.Range("F6").Value = 82

$wsMeans.Range("F7").Value = 6.4
$wsMeans.Range("G7").Value = 6

$wsMeans.Range("F8").Value = 4.3
$wsMeans.Range("G8").Value = 6.6

# --- Standard Deviations sheet ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

$wsSD.Range("G2").Value = 22

$wsSD.Range("F3").Value = 1.1
$wsSD.Range("G3").Value = 5.4

$wsSD.Range("F5").Value = 12
$wsSD.Range("G5").Value = 17

$wsSD.Range("F6").Value = 24
$wsSD.Range("G6").Value = 21

$wsSD.Range("F7").Value = 5.8

$wsSD.Range("F8").Value = 5.2
$wsSD.Range("G8").Value = 7

$wsSD.Range("F10").Value = 0.000000000000000016
